# Check rules & test updates
# Adds five new dictionary entries (aytk, waan, nabiip, nikdaa, algyax_) to the
# Sgüüxs dictionary sheet, and marks the "word" column (B) as Text-formatted
# so entries like these don't get auto-mangled by Excel's general formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Force the existing "word" column (B2:B91) to Text format -----------
# (mirrors selecting column B's data and setting Format Cells > Number > Text)
$ws.Range("B2:B91").NumberFormat = "@"

# --- 2. Row 92 (still inside the original table) ---------------------------
$ws.Range("B92").Value = "aytk"
$ws.Range("C92").Value = "name, identity"
$ws.Range("E92").Value = "noun"
$ws.Range("G92").Value = 1

# --- 3. New rows 93-96, extending the table ---------------------------------
$ws.Range("B93").Value = "'waan"
$ws.Range("C93").Value = "teeth"
$ws.Range("E93").Value = "noun"
$ws.Range("G93").Value = 1

$ws.Range("B94").Value = "nabiip"
$ws.Range("C94").Value = "uncle"
$ws.Range("E94").Value = "noun"
$ws.Range("G94").Value = 1

$ws.Range("B95").Value = "nikdaa"
$ws.Range("C95").Value = "aunt"
$ws.Range("E95").Value = "noun"
$ws.Range("G95").Value = 1

$ws.Range("B96").Value = "algyax_"
$ws.Range("C96").Value = "language, words, speech"
$ws.Range("E96").Value = "noun"
$ws.Range("G96").Value = 1

# --- 4. Extend the running Entry-ID counter in column A down to row 96 -----
$ws.Range("A93").Formula = "=1+A92"
$ws.Range("A94").Formula = "=1+A93"
$ws.Range("A95").Formula = "=1+A94"
$ws.Range("A96").Formula = "=1+A95"

# --- 5. Text-format the new word cells too ----------------------------------
$ws.Range("B92:B96").NumberFormat = "@"

# --- 6. B96 ("algyax_") picks up the same font treatment as the other
#        non-English-alphabet word entries (B86/B87, e.g. "'yaxwt") ---------
$ws.Range("B87").Copy()
$ws.Range("B96").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- 7. Leave the selection where the author ended up editing --------------
$ws.Range("C12").Select()

$wb.Save()
